$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("monsters")

# --- Add support for the new "Date" column (H) ---
$ws.Range("H1").Value = "Date"
$ws.Range("H2").Value = "date"
$ws.Range("H3").Value = 20190929
$ws.Range("H4").Value = 20190228

# Column D row 4 ("nums" field for the 2nd data record) is no longer populated
$ws.Range("D4").ClearContents()

# Give the new column a sensible width similar to the other data columns
$ws.Columns.Item(8).ColumnWidth = 11.4986979167
